$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

Set-TextValue 'D2' '40.012.15'
$ws.Range('E2').Value = '  +1.33%  '
Set-TextValue 'D3' '2.193.42'
$ws.Range('E3').Value = '  +1.29%  '
$ws.Range('E4').Value = '  +0.18%  '
Set-TextValue 'D5' '227.97'
$ws.Range('E5').Value = '  -0.54%  '
$ws.Range('E6').Value = '  +0.80%  '
Set-TextValue 'D7' '63.24'
$ws.Range('E7').Value = '  -0.67%  '
$ws.Range('E8').Value = '  +0.14%  '
Set-TextValue 'D9' '0.396'
$ws.Range('E9').Value = '  -0.64%  '
Set-TextValue 'D10' '0.0858'
$ws.Range('E10').Value = '  -0.91%  '
Set-TextValue 'D11' '0.103'
$ws.Range('E11').Value = '  -0.10%  '
Set-TextValue 'D12' '2.520.52'
$ws.Range('E12').Value = '  +1.44%  '
Set-TextValue 'D13' '15.76'
$ws.Range('E13').Value = '  -1.74%  '
$ws.Range('E14').Value = '  -1.13%  '
$ws.Range('E15').Value = '  -0.23%  '
$ws.Range('E16').Value = '  -0.46%  '
Set-TextValue 'D17' '2.180.65'
$ws.Range('E17').Value = '  +0.75%  '
Set-TextValue 'D18' '39.970.18'
$ws.Range('E18').Value = '  +1.28%  '
Set-TextValue 'D19' '0.0₃0904'
$ws.Range('E19').Value = '  +5.82%  '
Set-TextValue 'D20' '72.21'
$ws.Range('E20').Value = '  -0.03%  '
Set-TextValue 'D21' '6.06'
$ws.Range('E21').Value = '  -1.75%  '
Set-TextValue 'D22' '232.64'
$ws.Range('E22').Value = '  +1.48%  '
$ws.Range('E23').Value = '  +0.11%  '
Set-TextValue 'D24' '2.35'
$ws.Range('E24').Value = '  -0.73%  '
$ws.Range('E25').Value = '  +0.62%  '
$ws.Range('E26').Value = '  -1.15%  '
Set-TextValue 'D27' '171.87'
$ws.Range('E27').Value = '  -0.23%  '
Set-TextValue 'D28' '0.140'
$ws.Range('E28').Value = '  +2.13%  '
$ws.Range('E29').Value = '  +1.84%  '
Set-TextValue 'D30' '20.08'
$ws.Range('E30').Value = '  +1.77%  '
Set-TextValue 'D31' '2.73'
$ws.Range('E31').Value = '  +4.17%  '
$ws.Range('E32').Value = '  +0.26%  '
$ws.Range('E33').Value = '  -1.97%  '
$ws.Range('E34').Value = '  -2.67%  '
$ws.Range('E35').Value = '  -1.64%  '
$ws.Range('E36').Value = '  -0.16%  '
Set-TextValue 'D37' '3.86'
$ws.Range('E37').Value = '  +6.08%  '
$ws.Range('E38').Value = '  +0.01%  '
Set-TextValue 'D39' '5.01'
$ws.Range('E39').Value = '  +17.96%  '
$ws.Range('E40').Value = '  +0.15%  '
Set-TextValue 'D41' '102.88'
$ws.Range('E41').Value = '  -1.53%  '
$ws.Range('E42').Value = '  -1.02%  '
$ws.Range('B43').Value = 'InjectiveProtocol'
$ws.Range('C43').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
Set-TextValue 'D43' '17.44'
$ws.Range('E43').Value = '  -1.83%  '
$ws.Range('B44').Value = 'TrustWalletToken'
$ws.Range('C44').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
Set-TextValue 'D44' '1.23'
$ws.Range('E44').Value = '  +2.03%  '
$ws.Range('B45').Value = 'FraxShare'
$ws.Range('C45').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-TextValue 'D45' '8.31'
$ws.Range('E45').Value = '  +4.85%  '
$ws.Range('B46').Value = 'Maker'
$ws.Range('C46').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
Set-TextValue 'D46' '1.515.56'
$ws.Range('E46').Value = '  -1.52%  '
Set-TextValue 'D47' '0.0929'
$ws.Range('E47').Value = '  -0.56%  '
$ws.Range('E48').Value = '  -0.80%  '
$ws.Range('E49').Value = '  -0.40%  '
$ws.Range('E50').Value = '  +33.65%  '
Set-TextValue 'D51' '50.24'
$ws.Range('E51').Value = '  +7.51%  '
